$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PARENT_SITE_ID for the first two data rows was mistakenly set to 1413001.
# Correct it to the proper site id 1003001 (stored as a real number, not text).
$ws.Range("A2").Value = 1003001
$ws.Range("A3").Value = 1003001
